# Swap the bond detail columns (shortName, marketValue, couponRate,
# faceAmount, isin, percentWeight, cusip, sedol) between each pair of
# adjacent rows so the two "UST NOTE" / "US TREASURY N/B" entries that
# share the same maturity date trade places with each other.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns whose values get swapped between the paired rows.
$cols = @("D", "F", "G", "I", "K", "L", "Q", "R")

# These columns hold numeric-looking text (market value, coupon rate,
# face amount, percent weight) that must stay text -- otherwise Excel
# auto-converts them to numbers and we lose the original formatting
# (e.g. "4.250" -> 4.25, "0.90" -> 0.9).
$numericLookingCols = @("F", "G", "I", "L")

# Row pairs to swap (1-based, matching the worksheet row numbers).
$rowPairs = @(
    @(37, 38),
    @(41, 42),
    @(45, 46),
    @(47, 48),
    @(51, 52),
    @(55, 56)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $cell1 = $ws.Range("$col$r1")
        $cell2 = $ws.Range("$col$r2")

        # Use Value2 for reads (Value's getter round-trips oddly in this
        # runtime when chained directly into another Value assignment).
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        if ($numericLookingCols -contains $col) {
            # Force text so Excel doesn't silently reinterpret these as
            # numbers and normalize away the original text formatting.
            $cell1.NumberFormat = "@"
            $cell2.NumberFormat = "@"
        }

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
